# format_travel.xlsx: fix header labels (remove stray punctuation / use
# underscore-joined field names so they match the travel-admin migration),
# re-point the view to where the user left off editing, and force the
# sheet to print in portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row relabels (same columns, corrected text) ---------------
$ws.Range("A1").Value = "No"
$ws.Range("E1").Value = "Jml_Akreditasi"
$ws.Range("F1").Value = "tanggal_akreditasi"
$ws.Range("G1").Value = "lembaga_akreditasi"
$ws.Range("M1").Value = "kab_kota"

# --- Page setup: print portrait ----------------------------------------
$ws.PageSetup.Orientation = 1   # xlPortrait

# --- View / selection state ---------------------------------------------
# Leave the sheet scrolled/selected where the author last left it.
$ws.Range("M2").Select()
